$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LF / LF Lag coefficient table values are shrinking in magnitude / significance
# (Crisis and Credit Allocation controls added to the DAG-OLS regression).
# Column B = A, C = C, D = FFR, E = LF ; rows 2-5 = A/C/FFR/LF Lag.

$ws.Range("B2").Value = "-0.29***"
$ws.Range("B3").Value = "-1.02*"
$ws.Range("B4").Value = "0.06***"
$ws.Range("B5").Value = "0.01***"

$ws.Range("C2").Value = "-0.02***"
$ws.Range("C3").Value = "-0.42***"
$ws.Range("C4").Value = "-0.0*"
$ws.Range("C5").Value = "0.0***"

$ws.Range("D2").Value = "0.26*"
# "1.65" reads as a plain number to Excel, so force the cell to Text first
# (matches the source file, where this is stored as a shared text string).
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.65"
$ws.Range("D4").Value = "0.32***"
$ws.Range("D5").Value = "-0.01*"

$ws.Range("E2").Value = "-3.7*"
$ws.Range("E3:E5").NumberFormat = "@"
$ws.Range("E3").Value = "-1.49"
$ws.Range("E4").Value = "0.53"
$ws.Range("E5").Value = "0.1"
